$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new decision rows (7 and 8) describing the scene where
# the player is alone with Liz'Amar in Berto's room.
$ws.Range("C7").Value = "Con Liz'Amar solos en la habitación de Berto"
$ws.Range("D7").Value = "Pasar el rato con ella"
$ws.Range("F7").Value = "2"
$ws.Range("G7").Value = "-2"
$ws.Range("A7").Value = "6"
$ws.Range("B7").Value = "D"
$ws.Range("D8").Value = "Salir a entrenar"
$ws.Range("A8").Value = "7"
$ws.Range("B8").Value = "D"
$ws.Range("E7").Value = "4,-1"
$ws.Range("E8").Value = "4,-1"
$ws.Range("F8").Value = "-2"
$ws.Range("G8").Value = "2"

# Center the "name" column and merge it across the two new rows, matching
# the style used by the other multi-row decision blocks.
$ws.Range("C7:C8").HorizontalAlignment = -4108
$ws.Range("C7:C8").Merge()

# Leave the selection on E8, same as when the edit was made in the UI.
[void]$ws.Range("E8").Select()
